# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to stay a text/string value (matches the source data's
    # inline-string cells) instead of letting Excel auto-coerce plain
    # numeric-looking strings (e.g. "8.00") into a Number cell.
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

$ws.Range("D2").Value = "29.768.49"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.603.19"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +0.27%  "
Set-TextValue "D5" "212.57"
$ws.Range("E5").Value = "  -0.50%  "
Set-TextValue "D6" "0.518"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  +0.33%  "
Set-TextValue "D8" "28.53"
$ws.Range("E8").Value = "  +6.24%  "
$ws.Range("E9").Value = "  +2.02%  "
Set-TextValue "D10" "0.0604"
$ws.Range("E10").Value = "  +0.94%  "
Set-TextValue "D11" "0.0907"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "1.832.26"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "1.591.29"
$ws.Range("E13").Value = "  -2.59%  "
Set-TextValue "D14" "0.551"
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("D15").Value = "29.730.14"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("E16").Value = "  +0.66%  "
Set-TextValue "D17" "64.03"
$ws.Range("E17").Value = "  +0.81%  "
Set-TextValue "D18" "242.27"
$ws.Range("E18").Value = "  -0.99%  "
Set-TextValue "D19" "8.00"
$ws.Range("E19").Value = "  +4.73%  "
$ws.Range("D20").Value = "0.0₃0699"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("E21").Value = "  +0.32%  "
Set-TextValue "D22" "4.04"
$ws.Range("E22").Value = "  -0.42%  "
Set-TextValue "D23" "9.47"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("E24").Value = "  +1.18%  "
Set-TextValue "D25" "155.25"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("D34").Value = "1.422.72"
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("E35").Value = "  +3.86%  "
$ws.Range("E36").Value = "  -0.76%  "
Set-TextValue "D37" "2.86"
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  +1.57%  "
Set-TextValue "D40" "0.544"
$ws.Range("E40").Value = "  +1.80%  "
Set-TextValue "D41" "55.13"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  +5.65%  "
Set-TextValue "D43" "0.818"
$ws.Range("E43").Value = "  +2.41%  "
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "67.19"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D47" "0.995"
$ws.Range("E47").Value = "  +18.78%  "
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").Value = "1.741.87"
$ws.Range("E49").Value = "  -0.61%  "
Set-TextValue "D50" "86.61"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("E51").Value = "  +1.44%  "
